$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("for 10A")

# Re-measured data at the new voltage (235V) replaces rows 1-35;
# rows 36-61 (pre-existing lower-range measurements) are left untouched.
$ws.Cells.Item(1, 1).Value = 291
$ws.Cells.Item(1, 2).Value = 6930
$ws.Cells.Item(2, 1).Value = 289
$ws.Cells.Item(2, 2).Value = 6870
$ws.Cells.Item(3, 1).Value = 288
$ws.Cells.Item(3, 2).Value = 6870
$ws.Cells.Item(4, 1).Value = 355
$ws.Cells.Item(4, 2).Value = 8566
$ws.Cells.Item(5, 1).Value = 356
$ws.Cells.Item(5, 2).Value = 8600
$ws.Cells.Item(6, 1).Value = 353
$ws.Cells.Item(6, 2).Value = 8540
$ws.Cells.Item(7, 1).Value = 351
$ws.Cells.Item(7, 2).Value = 8490
$ws.Cells.Item(8, 1).Value = 351
$ws.Cells.Item(8, 2).Value = 8460
$ws.Cells.Item(9, 1).Value = 350
$ws.Cells.Item(9, 2).Value = 8439
$ws.Cells.Item(10, 1).Value = 344
$ws.Cells.Item(10, 2).Value = 8300
$ws.Cells.Item(11, 1).Value = 340
$ws.Cells.Item(11, 2).Value = 8250
$ws.Cells.Item(12, 1).Value = 340
$ws.Cells.Item(12, 2).Value = 8220
$ws.Cells.Item(13, 1).Value = 339
$ws.Cells.Item(13, 2).Value = 8170
$ws.Cells.Item(14, 1).Value = 333
$ws.Cells.Item(14, 2).Value = 8040
$ws.Cells.Item(15, 1).Value = 332
$ws.Cells.Item(15, 2).Value = 8000
$ws.Cells.Item(16, 1).Value = 331
$ws.Cells.Item(16, 2).Value = 7970
$ws.Cells.Item(17, 1).Value = 330
$ws.Cells.Item(17, 2).Value = 7947
$ws.Cells.Item(18, 1).Value = 329
$ws.Cells.Item(18, 2).Value = 7920
$ws.Cells.Item(19, 1).Value = 326
$ws.Cells.Item(19, 2).Value = 7880
$ws.Cells.Item(20, 1).Value = 321
$ws.Cells.Item(20, 2).Value = 7700
$ws.Cells.Item(21, 1).Value = 320
$ws.Cells.Item(21, 2).Value = 7676
$ws.Cells.Item(22, 1).Value = 312
$ws.Cells.Item(22, 2).Value = 7541
$ws.Cells.Item(23, 1).Value = 314
$ws.Cells.Item(23, 2).Value = 7522
$ws.Cells.Item(24, 1).Value = 311
$ws.Cells.Item(24, 2).Value = 7470
$ws.Cells.Item(25, 1).Value = 309
$ws.Cells.Item(25, 2).Value = 7426
$ws.Cells.Item(26, 1).Value = 308
$ws.Cells.Item(26, 2).Value = 7389
$ws.Cells.Item(27, 1).Value = 305
$ws.Cells.Item(27, 2).Value = 7350
$ws.Cells.Item(28, 1).Value = 293
$ws.Cells.Item(28, 2).Value = 7020
$ws.Cells.Item(29, 1).Value = 127
$ws.Cells.Item(29, 2).Value = 2645
$ws.Cells.Item(30, 1).Value = 126
$ws.Cells.Item(30, 2).Value = 2622
$ws.Cells.Item(31, 1).Value = 124
$ws.Cells.Item(31, 2).Value = 2588
$ws.Cells.Item(32, 1).Value = 122
$ws.Cells.Item(32, 2).Value = 2524
$ws.Cells.Item(33, 1).Value = 120
$ws.Cells.Item(33, 2).Value = 2455
$ws.Cells.Item(34, 1).Value = 119
$ws.Cells.Item(34, 2).Value = 2435
$ws.Cells.Item(35, 1).Value = 117
$ws.Cells.Item(35, 2).Value = 2390

# Add the Vrms/Irms conversion column (C = B/A) across the full data range.
$ws.Range("C1").Formula = "=B1/A1"
$ws.Range("C2:C33").Formula = "=B2/A2"
$ws.Range("C34:C61").Formula = "=B34/A34"

# Widen column C so the computed ratios are fully visible.
$ws.Columns.Item(3).ColumnWidth = 12.8888888888889

# Restore the on-screen selection/scroll position left by the edit.
$ws.Range("A36").Select()
